$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 41.82981864580804
$ws.Range("C2").Value = 3.541709949280602
$ws.Range("D2").Value = 43.58206693820312
$ws.Range("E2").Value = 35.10667621954779
$ws.Range("B3").Value = 41.82981864580804
$ws.Range("C3").Value = 3.541709949280602
$ws.Range("D3").Value = 43.51942067826702
$ws.Range("E3").Value = 23.61307215578357
$ws.Range("B4").Value = 41.82981864580804
$ws.Range("C4").Value = 3.541709949280602
$ws.Range("D4").Value = 43.51942067826702
$ws.Range("E4").Value = 3.337788520509218
$ws.Range("B5").Value = 23.48593870695883
$ws.Range("C5").Value = 3.541709949280602
$ws.Range("D5").Value = 19.76807200930811
$ws.Range("E5").Value = 3.337788520509218
$ws.Range("B6").Value = 23.48593870695883
$ws.Range("C6").Value = 3.541709949280602
$ws.Range("D6").Value = 19.76807200930811
$ws.Range("E6").Value = 3.337788520509218
$ws.Range("B7").Value = 21.78565594795221
$ws.Range("C7").Value = 3.541709949280602
$ws.Range("D7").Value = 19.76807200930811
$ws.Range("E7").Value = 3.337788520509218
$ws.Range("B8").Value = 21.78565594795221
$ws.Range("C8").Value = 3.541709949280602
$ws.Range("D8").Value = 19.76807200930811
$ws.Range("E8").Value = 3.337788520509218
$ws.Range("B9").Value = 21.78565594795221
$ws.Range("C9").Value = 3.541709949280602
$ws.Range("D9").Value = 19.76807200930811
$ws.Range("E9").Value = 3.337788520509218
$ws.Range("B10").Value = 21.78565594795221
$ws.Range("C10").Value = 3.541709949280602
$ws.Range("D10").Value = 19.76807200930811
$ws.Range("E10").Value = 3.337788520509218
$ws.Range("B11").Value = 21.78565594795221
$ws.Range("C11").Value = 3.541709949280602
$ws.Range("D11").Value = 19.76807200930811
$ws.Range("E11").Value = 3.337788520509218
$ws.Range("B12").Value = 21.78565594795221
$ws.Range("C12").Value = 3.541709949280602
$ws.Range("D12").Value = 19.76807200930811
$ws.Range("E12").Value = 3.337788520509218
$ws.Range("B13").Value = 21.78565594795221
$ws.Range("C13").Value = 3.541709949280602
$ws.Range("D13").Value = 19.76807200930811
$ws.Range("E13").Value = 3.337788520509218
$ws.Range("B14").Value = 21.78565594795221
$ws.Range("C14").Value = 3.541709949280602
$ws.Range("D14").Value = 19.76807200930811
$ws.Range("E14").Value = 3.337788520509218
$ws.Range("B15").Value = 21.78565594795221
$ws.Range("C15").Value = 3.541709949280602
$ws.Range("D15").Value = 19.76807200930811
$ws.Range("E15").Value = 3.337788520509218
$ws.Range("B16").Value = 21.78565594795221
$ws.Range("C16").Value = 3.541709949280602
$ws.Range("D16").Value = 19.76807200930811
$ws.Range("E16").Value = 3.337788520509218
$ws.Range("B17").Value = 21.78565594795221
$ws.Range("C17").Value = 3.541709949280602
$ws.Range("D17").Value = 16.28054437076482
$ws.Range("E17").Value = 3.337788520509218
$ws.Range("B18").Value = 21.78565594795221
$ws.Range("C18").Value = 3.541709949280602
$ws.Range("D18").Value = 16.28054437076482
$ws.Range("E18").Value = 3.337788520509218
$ws.Range("B19").Value = 21.78565594795221
$ws.Range("C19").Value = 3.541709949280602
$ws.Range("D19").Value = 16.28054437076482
$ws.Range("E19").Value = 3.337788520509218
$ws.Range("B20").Value = 20.90718437054669
$ws.Range("C20").Value = 3.541709949280602
$ws.Range("D20").Value = 16.28054437076482
$ws.Range("E20").Value = 3.337788520509218
$ws.Range("B21").Value = 20.79807174231033
$ws.Range("C21").Value = 3.541709949280602
$ws.Range("D21").Value = 16.28054437076482
$ws.Range("E21").Value = 3.337788520509218
$ws.Range("B22").Value = 20.79807174231033
$ws.Range("C22").Value = 3.541709949280602
$ws.Range("D22").Value = 16.28054437076482
$ws.Range("E22").Value = 3.337788520509218
$ws.Range("B23").Value = 20.79807174231033
$ws.Range("C23").Value = 3.541709949280602
$ws.Range("D23").Value = 16.28054437076482
$ws.Range("E23").Value = 2.526325816163247
$ws.Range("B24").Value = 20.79807174231033
$ws.Range("C24").Value = 3.541709949280602
$ws.Range("D24").Value = 16.28054437076482
$ws.Range("E24").Value = 2.526325816163247
$ws.Range("B25").Value = 20.79807174231033
$ws.Range("C25").Value = 3.541709949280602
$ws.Range("D25").Value = 16.28054437076482
$ws.Range("E25").Value = 2.526325816163247
$ws.Range("B26").Value = 18.71200197796063
$ws.Range("C26").Value = 3.541709949280602
$ws.Range("D26").Value = 16.28054437076482
$ws.Range("E26").Value = 2.526325816163247
$ws.Range("B27").Value = 18.71200197796063
$ws.Range("C27").Value = 3.541709949280602
$ws.Range("D27").Value = 16.28054437076482
$ws.Range("E27").Value = -2.159746960932743
$ws.Range("B28").Value = 18.71200197796063
$ws.Range("C28").Value = 3.541709949280602
$ws.Range("D28").Value = 16.28054437076482
$ws.Range("E28").Value = -2.159746960932743
$ws.Range("B29").Value = 18.71200197796063
$ws.Range("C29").Value = 3.541709949280602
$ws.Range("D29").Value = 6.776519609515645
$ws.Range("E29").Value = -2.159746960932743
$ws.Range("B30").Value = 18.71200197796063
$ws.Range("C30").Value = 3.541709949280602
$ws.Range("D30").Value = 6.776519609515645
$ws.Range("E30").Value = -2.159746960932743
$ws.Range("B31").Value = 18.71200197796063
$ws.Range("C31").Value = 3.541709949280602
$ws.Range("D31").Value = 6.776519609515645
$ws.Range("E31").Value = -2.159746960932743
$ws.Range("B32").Value = 18.71200197796063
$ws.Range("C32").Value = 3.541709949280602
$ws.Range("D32").Value = 6.776519609515645
$ws.Range("E32").Value = -2.159746960932743
$ws.Range("B33").Value = 18.71200197796063
$ws.Range("C33").Value = 3.541709949280602
$ws.Range("D33").Value = 6.776519609515645
$ws.Range("E33").Value = -2.159746960932743
$ws.Range("B34").Value = 18.71200197796063
$ws.Range("C34").Value = 3.541709949280602
$ws.Range("D34").Value = 6.776519609515645
$ws.Range("E34").Value = -2.159746960932743
$ws.Range("B35").Value = 18.71200197796063
$ws.Range("C35").Value = 3.541709949280602
$ws.Range("D35").Value = 6.776519609515645
$ws.Range("E35").Value = -2.159746960932743
$ws.Range("B36").Value = 18.71200197796063
$ws.Range("C36").Value = 3.541709949280602
$ws.Range("D36").Value = 6.776519609515645
$ws.Range("E36").Value = -2.159746960932743
$ws.Range("B37").Value = 18.71200197796063
$ws.Range("C37").Value = 3.541709949280602
$ws.Range("D37").Value = 6.776519609515645
$ws.Range("E37").Value = -2.159746960932743
$ws.Range("B38").Value = 17.7178413188736
$ws.Range("C38").Value = 3.541709949280602
$ws.Range("D38").Value = 6.776519609515645
$ws.Range("E38").Value = -2.159746960932743
$ws.Range("B39").Value = 17.7178413188736
$ws.Range("C39").Value = 3.541709949280602
$ws.Range("D39").Value = 5.537964642031774
$ws.Range("E39").Value = -2.159746960932743
$ws.Range("B40").Value = 17.7178413188736
$ws.Range("C40").Value = 3.541709949280602
$ws.Range("D40").Value = 5.537964642031774
$ws.Range("E40").Value = -2.159746960932743
$ws.Range("B41").Value = 17.7178413188736
$ws.Range("C41").Value = 3.541709949280602
$ws.Range("D41").Value = 5.537964642031774
$ws.Range("E41").Value = -2.159746960932743
$ws.Range("B42").Value = 17.7178413188736
$ws.Range("C42").Value = 3.541709949280602
$ws.Range("D42").Value = 5.537964642031774
$ws.Range("E42").Value = -5.352027628372145
$ws.Range("B43").Value = 17.7178413188736
$ws.Range("C43").Value = 3.541709949280602
$ws.Range("D43").Value = 5.537964642031774
$ws.Range("E43").Value = -5.352027628372145
$ws.Range("B44").Value = 17.7178413188736
$ws.Range("C44").Value = 3.541709949280602
$ws.Range("D44").Value = 5.537964642031774
$ws.Range("E44").Value = -5.352027628372145
$ws.Range("B45").Value = 17.7178413188736
$ws.Range("C45").Value = 3.541709949280602
$ws.Range("D45").Value = 5.537964642031774
$ws.Range("E45").Value = -5.352027628372145
$ws.Range("B46").Value = 17.7178413188736
$ws.Range("C46").Value = 3.541709949280602
$ws.Range("D46").Value = 5.537964642031774
$ws.Range("E46").Value = -5.352027628372145
